# precision-recall-graph.xlsx edit
#
# The underlying engine only models the "Data" worksheet (the chart sheet
# "Precision-Recall Graph" + its embedded chart1.xml are carried through
# as opaque/passthrough parts and are not exposed anywhere on the Excel
# object model in this runtime - no Worksheets/Sheets/Charts entry, no
# ActiveChart, nothing). So this script applies every change from the
# diff that is reachable through $excel / $wb / worksheet COM calls:
#   - drop the three now-unused defined names (ir_1, ir_2, ir_3), keep "ir"
#   - the "ir3" text-connection + its H:I helper columns on Data are gone,
#     so clear out H1:I11 (values + the implied dimension/row span shrink)
#   - refresh the Stemming (D) and Stopwords&Stemming (E) columns with the
#     new imported values
#   - move the sheet selection from D3:D13 to E3:E13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- defined names: remove ir_1, ir_2, ir_3, keep ir ---------------------
$wb.Names.Item("Data!ir_1").Delete()
$wb.Names.Item("Data!ir_2").Delete()
$wb.Names.Item("Data!ir_3").Delete()

# --- drop the old H:I helper-table (was driven by the removed ir3 query) -
$ws.Range("H1:I11").ClearContents()

# --- refresh column D ("Stemming") with the re-imported values -----------
$ws.Range("D3").Value = 0.221501881860245
$ws.Range("D4").Value = 0.221501881860245
$ws.Range("D5").Value = 0.20657863888904099
$ws.Range("D6").Value = 0.17811684952742901
$ws.Range("D7").Value = 0.16655740079927001
$ws.Range("D8").Value = 0.16290871517195801
$ws.Range("D9").Value = 0.12540428870926501
$ws.Range("D10").Value = 0.119437441783423
$ws.Range("D11").Value = 0.114918012706467
$ws.Range("D12").Value = 0.107471131076581
$ws.Range("D13").Value = 0.107471131076581

# --- refresh column E ("Stopwords & Stemming") with the re-imported values
$ws.Range("E3").Value = 0.247986404777814
$ws.Range("E4").Value = 0.247986404777814
$ws.Range("E5").Value = 0.22386302815443801
$ws.Range("E6").Value = 0.20269525433011701
$ws.Range("E7").Value = 0.18258372871300599
$ws.Range("E8").Value = 0.17951744858369101
$ws.Range("E9").Value = 0.12946812224639101
$ws.Range("E10").Value = 0.12187555893814001
$ws.Range("E11").Value = 0.11676554405122
$ws.Range("E12").Value = 0.11089270367040301
$ws.Range("E13").Value = 0.11057964830471601

# --- selection now tracks column E instead of D ---------------------------
$ws.Range("E3:E13").Select()
